# Apply "fixed param talynM txbrdfw" edit:
#  - Column P ("note") is repurposed to hold "PacketLengthPayload" values
#    (mirrors the SilenceDuration column O).
#  - A new column Q ("tx_gain_row") also mirrors SilenceDuration (col O).
#  - A new column R takes over as the "note" column, carrying the text
#    that used to live in column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("P1").Value = "PacketLengthPayload"
$ws.Range("Q1").Value = "tx_gain_row"
$ws.Range("R1").Value = "note"

# --- Data rows (2 through 17) ---
$lastRow = 17
for ($r = 2; $r -le $lastRow; $r++) {
    $silenceDuration = $ws.Cells.Item($r, 15).Value2   # column O
    $oldNote = $ws.Cells.Item($r, 16).Value2            # column P (old "note")

    $ws.Cells.Item($r, 16).Value = $silenceDuration     # column P -> PacketLengthPayload
    $ws.Cells.Item($r, 17).Value = $silenceDuration     # column Q -> tx_gain_row
    $ws.Cells.Item($r, 18).Value = $oldNote             # column R -> note (moved)
}
